# "Started working on Read-in of region selection"
#
# The helper columns G:I (Country-Node / Country_Sublevel_1_Node /
# Country_Sublevel_2_Node, with their 0/1 "node" flags) are removed from
# the Region_selection sheet, and a few of the remaining "Region Selected"
# flags are reset from 1 to 0.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Region_selection")

# Remove the now-unused helper columns G, H and I (and their header/shared
# strings) entirely - this also shrinks the sheet dimension back to A1:E10.
$ws1.Range("G1:I10").Delete()

# A few rows are no longer flagged as a selected region.
$ws1.Range("E5").Value2 = 0
$ws1.Range("E7").Value2 = 0
$ws1.Range("E8").Value2 = 0

# Leave the cursor where the author left it while starting to work on the
# region-selection read-in.
$ws1.Range("G23").Select()
